$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Acc_num" (DE2) description cell (row 3, column C): append the new
# "It is used to review transactions" sentence to the existing text.
$ws.Range("C3").Value = "Account number of the customer will be used to deposit,withdrawal, verififcation, and other banking activities. It is used to review transactions"

# The longer text now wraps onto two lines, so the row grows to match
# the height Excel would auto-compute for the wrapped text (same height
# used elsewhere in the sheet for two-line wrapped descriptions).
$ws.Rows.Item(3).RowHeight = 31.5

# Reflect where the user's cursor/selection ended up after making the edit.
$ws.Range("C3").Select()
